# Scheduled-runner price refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H-N) for the leves whose backing market data moved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 375.92426
$ws.Range("J17").Value = 349.36066
$ws.Range("L17").Value = 1048.08198
$ws.Range("N17").Value = -1384.08198

$ws.Range("H100").Value = 18520864
$ws.Range("I100").Value = 30304586
$ws.Range("K100").Value = 30304586
$ws.Range("M100").Value = -30304045

$ws.Range("H106").Value = 1932.6364
$ws.Range("I106").Value = 1932.6364
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 1932.6364
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -1301.6364
$ws.Range("N106").ClearContents()

$ws.Range("H129").Value = 937.8570999999999
$ws.Range("I129").Value = 282.875
$ws.Range("J129").Value = 1199.85
$ws.Range("K129").Value = 848.625
$ws.Range("L129").Value = 3599.55
$ws.Range("M129").Value = 4151.375
$ws.Range("N129").Value = -13599.55

$ws.Range("H137").Value = 908.0789
$ws.Range("I137").Value = 808.48
$ws.Range("K137").Value = 2425.44
$ws.Range("M137").Value = 124.5599999999999

$ws.Range("H138").Value = 3399.1194
$ws.Range("I138").Value = 2104.0435
$ws.Range("J138").Value = 4076.0908
$ws.Range("K138").Value = 6312.130500000001
$ws.Range("L138").Value = 12228.2724
$ws.Range("M138").Value = -1172.130500000001
$ws.Range("N138").Value = -22508.2724


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1243.3684
$ws.Range("I61").Value = 1379.1111
$ws.Range("J61").Value = 1121.2
$ws.Range("K61").Value = 1379.1111
$ws.Range("L61").Value = 1121.2
$ws.Range("M61").Value = -1167.1111
$ws.Range("N61").Value = -1545.2

$ws.Range("H74").Value = 1042.7028
$ws.Range("I74").Value = 962.73334
$ws.Range("J74").Value = 1385.4286
$ws.Range("K74").Value = 962.73334
$ws.Range("L74").Value = 1385.4286
$ws.Range("M74").Value = -88.73334
$ws.Range("N74").Value = -3133.4286

$ws.Range("H77").Value = 1042.7028
$ws.Range("I77").Value = 962.73334
$ws.Range("J77").Value = 1385.4286
$ws.Range("K77").Value = 4813.6667
$ws.Range("L77").Value = 6927.143
$ws.Range("M77").Value = -445.6666999999998
$ws.Range("N77").Value = -15663.143

$ws.Range("H88").Value = 4583.9165
$ws.Range("I88").Value = 2600
$ws.Range("J88").Value = 6001
$ws.Range("K88").Value = 2600
$ws.Range("L88").Value = 6001
$ws.Range("M88").Value = -2194
$ws.Range("N88").Value = -6813

$ws.Range("H91").Value = 4583.9165
$ws.Range("I91").Value = 2600
$ws.Range("J91").Value = 6001
$ws.Range("K91").Value = 2600
$ws.Range("L91").Value = 6001
$ws.Range("M91").Value = -1196
$ws.Range("N91").Value = -8809

$ws.Range("H132").Value = 25026922
$ws.Range("I132").Value = 33335070
$ws.Range("J132").Value = 102474.8
$ws.Range("K132").Value = 100005210
$ws.Range("L132").Value = 307424.4
$ws.Range("M132").Value = -100002680
$ws.Range("N132").Value = -312484.4

$ws.Range("H136").Value = 1243.3684
$ws.Range("I136").Value = 1379.1111
$ws.Range("J136").Value = 1121.2
$ws.Range("K136").Value = 4137.3333
$ws.Range("L136").Value = 3363.6
$ws.Range("M136").Value = -1587.3333
$ws.Range("N136").Value = -8463.6


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 24205.625
$ws.Range("I134").Value = 9009.200000000001
$ws.Range("J134").Value = 49533
$ws.Range("K134").Value = 27027.6
$ws.Range("L134").Value = 148599
$ws.Range("M134").Value = -24492.6
$ws.Range("N134").Value = -153669


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9617880
$ws.Range("J31").Value = 5450
$ws.Range("L31").Value = 5450
$ws.Range("N31").Value = -6040

$ws.Range("H34").Value = 9617880
$ws.Range("J34").Value = 5450
$ws.Range("L34").Value = 5450
$ws.Range("N34").Value = -5854

$ws.Range("H58").Value = 903.6875
$ws.Range("I58").Value = 935.3077
$ws.Range("J58").Value = 766.6667
$ws.Range("K58").Value = 935.3077
$ws.Range("L58").Value = 766.6667
$ws.Range("M58").Value = -732.3077
$ws.Range("N58").Value = -1172.6667

$ws.Range("H105").Value = 1091.4706
$ws.Range("I105").Value = 860.4545000000001
$ws.Range("J105").Value = 1515
$ws.Range("K105").Value = 860.4545000000001
$ws.Range("L105").Value = 1515
$ws.Range("M105").Value = 886.5454999999999
$ws.Range("N105").Value = -5009

$ws.Range("H132").Value = 58794.277
$ws.Range("I132").Value = 2191.5833
$ws.Range("J132").Value = 171999.67
$ws.Range("K132").Value = 6574.749899999999
$ws.Range("L132").Value = 515999.01
$ws.Range("M132").Value = -4044.749899999999
$ws.Range("N132").Value = -521059.01

$ws.Range("H134").Value = 2961.125
$ws.Range("I134").Value = 1877
$ws.Range("J134").Value = 10550
$ws.Range("K134").Value = 5631
$ws.Range("L134").Value = 31650
$ws.Range("M134").Value = -3096
$ws.Range("N134").Value = -36720

$ws.Range("H136").Value = 903.6875
$ws.Range("I136").Value = 935.3077
$ws.Range("J136").Value = 766.6667
$ws.Range("K136").Value = 2805.9231
$ws.Range("L136").Value = 2300.0001
$ws.Range("M136").Value = -255.9231
$ws.Range("N136").Value = -7400.0001


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 911023.5600000001
$ws.Range("I132").Value = 834422
$ws.Range("J132").Value = 1002945.4
$ws.Range("K132").Value = 2503266
$ws.Range("L132").Value = 3008836.2
$ws.Range("M132").Value = -2500736
$ws.Range("N132").Value = -3013896.2


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2711.7
$ws.Range("I7").Value = 1249.8889
$ws.Range("J7").Value = 3907.7273
$ws.Range("K7").Value = 1249.8889
$ws.Range("L7").Value = 3907.7273
$ws.Range("M7").Value = -1137.8889
$ws.Range("N7").Value = -4131.7273

$ws.Range("H40").Value = 1991.4546
$ws.Range("I40").Value = 1790.8182
$ws.Range("K40").Value = 1790.8182
$ws.Range("M40").Value = -1654.8182

$ws.Range("H126").Value = 2711.7
$ws.Range("I126").Value = 1249.8889
$ws.Range("J126").Value = 3907.7273
$ws.Range("K126").Value = 3749.6667
$ws.Range("L126").Value = 11723.1819
$ws.Range("M126").Value = -1279.6667
$ws.Range("N126").Value = -16663.1819

$ws.Range("H132").Value = 295902.53
$ws.Range("I132").Value = 371842.84
$ws.Range("J132").Value = 2989.8572
$ws.Range("K132").Value = 1115528.52
$ws.Range("L132").Value = 8969.571599999999
$ws.Range("M132").Value = -1112998.52
$ws.Range("N132").Value = -14029.5716

$ws.Range("H136").Value = 5450.5674
$ws.Range("I136").Value = 5568.5386
$ws.Range("J136").Value = 5171.727
$ws.Range("K136").Value = 16705.6158
$ws.Range("L136").Value = 15515.181
$ws.Range("M136").Value = -14155.6158
$ws.Range("N136").Value = -20615.181


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 52098436
$ws.Range("I132").Value = 66471804
$ws.Range("J132").Value = 3228985.2
$ws.Range("K132").Value = 199415412
$ws.Range("L132").Value = 9686955.600000001
$ws.Range("M132").Value = -199412882
$ws.Range("N132").Value = -9692015.600000001

$ws.Range("H136").Value = 32293.938
$ws.Range("I136").Value = 39450.19
$ws.Range("J136").Value = 1283.5
$ws.Range("K136").Value = 118350.57
$ws.Range("L136").Value = 3850.5
$ws.Range("M136").Value = -115800.57
$ws.Range("N136").Value = -8950.5
